$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the label in G1 from "Adicion" to "Añadicion"
$ws.Range("G1").Value = "Añadicion"

# Move the selection/active cell to G5
$ws.Activate()
$ws.Range("G5").Select()
